$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the B-column figures (now fractions formatted as percentages) ---
$ws.Range("B2").Value = -1.1863039602580059
$ws.Range("B3").Value = -0.87730684124276759
$ws.Range("B4").Value = -0.98865385341600798
$ws.Range("B5").Value = -0.27508513225814446
$ws.Range("B6").Value = -0.151788671431939
$ws.Range("B7").Value = -0.13648917617817061
$ws.Range("B8").Value = 0.0099620728755395495
$ws.Range("B9").Value = 0.61804224730370816
$ws.Range("B10").Value = 0.76268770220451798

# C4 count bumped
$ws.Range("C4").Value = 15

# Format the B column figures as a one-decimal percentage
$ws.Range("B2:B10").NumberFormat = "0.0%"

# Give the top two region labels a distinguishing look: black Calibri text,
# a thin light-gray box around the left/right edges, and wrapped text.
$labels = $ws.Range("A2:A3")
$labels.Borders.Item(7).LineStyle = 1
$labels.Borders.Item(7).Color = 12632256
$labels.Borders.Item(10).LineStyle = 1
$labels.Borders.Item(10).Color = 12632256
$labels.Font.Color = 0
$labels.WrapText = $true

# Move the active selection down to A2:C10 with A2 as the active cell
[void]$ws.Range("A2:C10").Select()
